$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Entity1 only had one attribute row (Attr1, row 4) but should have all three
# (Attr1, Attr2, Attr3), matching the other entities. Insert two rows right
# after the existing Attr1 row and fill them in with the missing attributes,
# copying the formatting that is already used for the other "Attr" rows.
$ws.Rows.Item(5).Resize(2, 1).EntireRow.Insert()

$ws.Range("A5").Value = "Attr2"
$ws.Range("A6").Value = "Attr3"

# Entity2's value was wrong (2 instead of 3) - fix it now that the row has
# shifted down to row 7.
$ws.Range("B7").Value = 3

# Reflect the last place the user was working.
$ws.Range("D10").Select()
